$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Helper: force a run split at a given zero-length character position by
# dropping a temporary bookmark there and immediately deleting it again.
# (Word's own run-splitting logic kicks in while the bookmark is anchored
# at that offset, and the resulting run boundary survives the bookmark's
# removal.)
# ----------------------------------------------------------------------
function Split-RunAt($pos, $bmName) {
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($bmName, $r)
    $bm = $d.Bookmarks.Item($bmName)
    $bm.Delete()
}

# ------------------------------------------------------------------
# 1) Title paragraph: "cenarij 4: " -> "Scenarij 4: "
#    Inserted as its own leading run ("S"), not merged into the run
#    that already holds "cenarij 4: ".
# ------------------------------------------------------------------
$titleStart = $d.Range(0, 0)
$titleStart.InsertBefore("S")
Split-RunAt 1 "_TMP_TITLE"

# ------------------------------------------------------------------
# 2) Remove the stray "_GoBack" bookmark that currently sits, on its
#    own, in the empty trailing paragraph of the document. It gets
#    re-created (below) further up in the document, where the cursor
#    genuinely was the last time the file was saved.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 3) "Korisnik se nalazi na kućnoj adresi, korisnik posjeduje novac da
#    plati uređaj" — text itself is unchanged, but a "_GoBack" bookmark
#    now sits between "uređa" and the final "j" (this is where editing
#    last left off).
# ------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("posjeduje novac da plati uređaj")
$bookmarkPos = $r.End - 1
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 4) Typo fix: "urešaj" -> "uređaj" in "4. Dostavljač predaje urešaj
#    korisniku", ending up split as two runs ("...uređ" / "aj
#    korisniku").
# ------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("urešaj korisniku", $true, $false, $false, $false, $false, $true, 1, $false, "uređaj korisniku", 2)

$r3 = $d.Content
$null = $r3.Find.Execute("Dostavljač predaje uređ")
Split-RunAt $r3.End "_TMP_UREDAJ"
